# Auto-update draw results: append the 2025-11-12 Pick 3 draw as a new
# row (row 57) at the bottom of the Results sheet, mirroring every
# existing row exactly (all five columns stored as plain text, not as
# numbers/dates that happen to look similar).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 57

# Force the new row to be treated as plain text before writing any
# values into it, so e.g. "2025-11-12" and "251112" stay text instead of
# being auto-converted to a date serial / number (matches every other
# row already in the sheet).
$rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5))
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-11-12"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251112"
$ws.Cells.Item($newRow, 4).Value = "8-9-5"
$ws.Cells.Item($newRow, 5).Value = "2025-11-12T21:40:14.496+04:00"
